$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 151991.61
$ws.Range("I6").Value = 200182.92
$ws.Range("K6").Value = 600548.76
$ws.Range("M6").Value = -600436.76
$ws.Range("H51").Value = 6040.2
$ws.Range("J51").Value = 6578.6
$ws.Range("L51").Value = 6578.6
$ws.Range("N51").Value = -7546.6
$ws.Range("H125").Value = 5884852.5
$ws.Range("I125").Value = 5000
$ws.Range("K125").Value = 45000
$ws.Range("M125").Value = -42540
$ws.Range("H135").Value = 888.0357
$ws.Range("I135").Value = 784.2381
$ws.Range("J135").Value = 1199.4286
$ws.Range("K135").Value = 7058.142900000001
$ws.Range("L135").Value = 10794.8574
$ws.Range("M135").Value = -4523.142900000001
$ws.Range("N135").Value = -15864.8574

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2491.85
$ws.Range("I61").Value = 2322.4666
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 2322.4666
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -2110.4666
$ws.Range("N61").Value = -3424
$ws.Range("H74").Value = 8708.823
$ws.Range("I74").Value = 15342
$ws.Range("J74").Value = 2812.6667
$ws.Range("K74").Value = 15342
$ws.Range("L74").Value = 2812.6667
$ws.Range("M74").Value = -14468
$ws.Range("N74").Value = -4560.6667
$ws.Range("H77").Value = 8708.823
$ws.Range("I77").Value = 15342
$ws.Range("J77").Value = 2812.6667
$ws.Range("K77").Value = 76710
$ws.Range("L77").Value = 14063.3335
$ws.Range("M77").Value = -72342
$ws.Range("N77").Value = -22799.3335
$ws.Range("H110").Value = 125390.5
$ws.Range("I110").Value = 167083
$ws.Range("J110").Value = 313
$ws.Range("K110").Value = 167083
$ws.Range("L110").Value = 313
$ws.Range("M110").Value = -165038
$ws.Range("N110").Value = -4403
$ws.Range("H125").Value = 34905
$ws.Range("J125").Value = 34905
$ws.Range("L125").Value = 34905
$ws.Range("N125").Value = -44745
$ws.Range("H132").Value = 2134.1538
$ws.Range("I132").Value = 1480.3572
$ws.Range("J132").Value = 2896.9167
$ws.Range("K132").Value = 4441.071599999999
$ws.Range("L132").Value = 8690.750100000001
$ws.Range("M132").Value = -1911.071599999999
$ws.Range("N132").Value = -13750.7501
$ws.Range("H136").Value = 2491.85
$ws.Range("I136").Value = 2322.4666
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 6967.399800000001
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -4417.399800000001
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 33754
$ws.Range("J51").Value = 33754
$ws.Range("L51").Value = 33754
$ws.Range("N51").Value = -34736
$ws.Range("H97").Value = 2000
$ws.Range("I97").Value = 2000
$ws.Range("K97").Value = 2000
$ws.Range("M97").Value = -1009
$ws.Range("H105").Value = 2791.7058
$ws.Range("I105").Value = 2713.5
$ws.Range("K105").Value = 2713.5
$ws.Range("M105").Value = -966.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2511.2693
$ws.Range("I58").Value = 2292.2856
$ws.Range("J58").Value = 2766.75
$ws.Range("K58").Value = 2292.2856
$ws.Range("L58").Value = 2766.75
$ws.Range("M58").Value = -2089.2856
$ws.Range("N58").Value = -3172.75
$ws.Range("H134").Value = 2628.889
$ws.Range("I134").Value = 2887.0667
$ws.Range("J134").Value = 1338
$ws.Range("K134").Value = 8661.2001
$ws.Range("L134").Value = 4014
$ws.Range("M134").Value = -6126.2001
$ws.Range("N134").Value = -9084
$ws.Range("H136").Value = 2511.2693
$ws.Range("I136").Value = 2292.2856
$ws.Range("J136").Value = 2766.75
$ws.Range("K136").Value = 6876.8568
$ws.Range("L136").Value = 8300.25
$ws.Range("M136").Value = -4326.8568
$ws.Range("N136").Value = -13400.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = ""
$ws.Range("N7").Value = 0
$ws.Range("H80").Value = 4980
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 4980
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = ""
$ws.Range("M80").Value = 14940
$ws.Range("N80").Value = -16812
$ws.Range("H83").Value = 4980
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 4980
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = ""
$ws.Range("M83").Value = 44820
$ws.Range("N83").Value = -54180
$ws.Range("H92").Value = 818.1429000000001
$ws.Range("I92").Value = 1051
$ws.Range("J92").Value = 725
$ws.Range("K92").Value = 3153
$ws.Range("L92").Value = 2175
$ws.Range("M92").Value = -1905
$ws.Range("N92").Value = -4671
$ws.Range("H131").Value = 854.3099999999999
$ws.Range("J131").Value = 857.8889
$ws.Range("L131").Value = 2573.6667
$ws.Range("N131").Value = -12653.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5049.8726
$ws.Range("I70").Value = 4454.793
$ws.Range("J70").Value = 5713.615
$ws.Range("K70").Value = 4454.793
$ws.Range("L70").Value = 5713.615
$ws.Range("M70").Value = -4184.793
$ws.Range("N70").Value = -6253.615
$ws.Range("H73").Value = 5049.8726
$ws.Range("I73").Value = 4454.793
$ws.Range("J73").Value = 5713.615
$ws.Range("K73").Value = 4454.793
$ws.Range("L73").Value = 5713.615
$ws.Range("M73").Value = -3518.793
$ws.Range("N73").Value = -7585.615
$ws.Range("H80").Value = 2335.3333
$ws.Range("I80").Value = 2000
$ws.Range("J80").Value = 2503
$ws.Range("K80").Value = 2000
$ws.Range("L80").Value = 2503
$ws.Range("M80").Value = -1002
$ws.Range("N80").Value = -4499
$ws.Range("H83").Value = 2335.3333
$ws.Range("I83").Value = 2000
$ws.Range("J83").Value = 2503
$ws.Range("K83").Value = 10000
$ws.Range("L83").Value = 12515
$ws.Range("M83").Value = -5008
$ws.Range("N83").Value = -22499

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 27400
$ws.Range("J14").Value = 27400
$ws.Range("L14").Value = 27400
$ws.Range("N14").Value = -27744
$ws.Range("H82").Value = 2111.158
$ws.Range("I82").Value = 1812.75
$ws.Range("J82").Value = 2328.182
$ws.Range("K82").Value = 1812.75
$ws.Range("L82").Value = 2328.182
$ws.Range("M82").Value = -1451.75
$ws.Range("N82").Value = -3050.182
$ws.Range("H85").Value = 2111.158
$ws.Range("I85").Value = 1812.75
$ws.Range("J85").Value = 2328.182
$ws.Range("K85").Value = 1812.75
$ws.Range("L85").Value = 2328.182
$ws.Range("M85").Value = -564.75
$ws.Range("N85").Value = -4824.182
$ws.Range("H103").Value = 20167.334
$ws.Range("J103").Value = 20167.334
$ws.Range("L103").Value = 20167.334
$ws.Range("N103").Value = -22511.334
$ws.Range("H133").Value = 32179.875
$ws.Range("J133").Value = 32179.875
$ws.Range("L133").Value = 32179.875
$ws.Range("N133").Value = -37239.875
$ws.Range("H136").Value = 20837032
$ws.Range("I136").Value = 3722.9167
$ws.Range("J136").Value = 83336960
$ws.Range("K136").Value = 11168.7501
$ws.Range("L136").Value = 250010880
$ws.Range("M136").Value = -8618.750100000001
$ws.Range("N136").Value = -250015980

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 45000.8
$ws.Range("I5").Value = 5000
$ws.Range("J5").Value = 55001
$ws.Range("K5").Value = 5000
$ws.Range("L5").Value = 55001
$ws.Range("M5").Value = -4888
$ws.Range("N5").Value = -55225
$ws.Range("H96").Value = 1196
$ws.Range("I96").Value = 726.6667
$ws.Range("K96").Value = 726.6667
$ws.Range("M96").Value = 646.3333
$ws.Range("H100").Value = 14397.25
$ws.Range("I100").Value = 695.3333
$ws.Range("J100").Value = 55503
$ws.Range("K100").Value = 1390.6666
$ws.Range("L100").Value = 111006
$ws.Range("M100").Value = -849.6666
$ws.Range("N100").Value = -112088
$ws.Range("H126").Value = 144514.14
$ws.Range("I126").Value = 200938.8
$ws.Range("K126").Value = 602816.3999999999
$ws.Range("M126").Value = -600346.3999999999
